$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Review")

# --- Row 5: add sequence number ---
$ws.Range("A5").Value = 4

# --- Row 6: add sequence number, fix file reference text ---
$ws.Range("A6").Value = 5

# --- New shared strings must be created in this exact order so the ---
# --- sharedStrings table matches the authored edit order.          ---
$ws.Range("I7").Value = "Tutorial_07 `n1) Remove unnecessary file"
$ws.Range("F6").Value = "Tutorial_06`nindex.php"
$ws.Range("F7").Value = "Tutorial_07`nindex.php"
$ws.Range("F7").WrapText = $true
$ws.Range("P7").Value = "1)removed unnecessary file"

# --- Row 7: fill in the rest of the tutorial-7 review entry ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 44552
$ws.Range("C7").Value = "Open"
$ws.Range("D7").Value = "Others"
$ws.Range("E7").Value = "Improvement"
$ws.Range("O7").Value = "PyaePyaeHan"
$ws.Range("V7").Value = 44552
$ws.Range("V7").NumberFormat = "mm-dd-yy"
$ws.Range("W7").Value = "EiWahWin"
$ws.Range("X7").Value = "Done"

# --- Update the view: scroll back to the top-left and reselect X7 ---
$ws.Range("X7").Select()

# --- Define the print area for the Review sheet ---
$ws.PageSetup.PrintArea = '$A$1:$X$18'
